$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Values for column I (I0) and column J (IF), rows 2 through 61
$iValues = @(8, 4, 6, 6, 6, 9, 7, 6, 4, 3, 9, 3, 7, 5, 8, 5, 4, 9, 5, 8, 8, 8, 5, 8, 6, 6, 9, 7, 6, 8, 7, 9, 7, 9, 7, 9, 7, 4, 8, 9, 8, 6, 3, 4, 4, 7, 7, 8, 9, 10, 8, 4, 7, 9, 3, 9, 5, 7, 5, 7)
$jValues = @(8, 4, 6, 7, 7, 9, 7, 7, 5, 5, 9, 4, 7, 6, 8, 6, 6, 9, 5, 8, 8, 8, 5, 9, 6, 7, 9, 7, 6, 8, 7, 9, 7, 9, 8, 9, 7, 6, 9, 9, 8, 7, 5, 6, 5, 7, 7, 8, 9, 10, 8, 6, 8, 9, 3, 9, 5, 7, 5, 7)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
